$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.5316996666666666
$ws.Range("M2").Value = 14.11187666666667
$ws.Range("N2").Value = 42.33562999999999
$ws.Range("O2").Value = 0.08862966207485527
$ws.Range("P2").Value = 0.08862966207485526
$ws.Range("Q2").Value = 7.503280119707775
$ws.Range("R2").Value = 67.52952107736998
$ws.Range("S2").Value = 0.08862966207485527
$ws.Range("T2").Value = 0.08862966207485526

# Row 3
$ws.Range("G3").Value = 0.5316996666666666
$ws.Range("O3").Value = 0.7176943460983047
$ws.Range("P3").Value = 0.7176943460983046
$ws.Range("Q3").Value = 60.75913630989521
$ws.Range("R3").Value = 546.8322267890569
$ws.Range("S3").Value = 0.7176943460983047
$ws.Range("T3").Value = 0.7176943460983046

# Row 4
$ws.Range("G4").Value = 0.5316996666666666
$ws.Range("O4").Value = 0.1936759918268401
$ws.Range("P4").Value = 0.1936759918268401
$ws.Range("Q4").Value = 16.39637549234544
$ws.Range("S4").Value = 0.1936759918268401
$ws.Range("T4").Value = 0.1936759918268401
